$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend rows 13-17 by copying the formatting (and border style on column A)
# from the last existing row (row 12), so the new rows inherit the same look.
for ($r = 13; $r -le 17; $r++) {
    $ws.Range("A12:K12").Copy($ws.Range("A" + $r + ":K" + $r))
}

# Final values for rows 2-17, columns A-K (HOURS table data window shifted + 5 new rows)
$rowData = @(
    @(29, 30, 98.71785172641819, 0, 0, 0, 148.0767775896273, 0.37, 0.46, 0.15, 0.33),
    @(30, 31, 130.7873219730109, 0, 0, 0, 196.1809829595164, 0.37, 0.49, 0.13, 0.4),
    @(31, 32, 201.8420447128433, 0, 0, 0, 302.7630670692649, 0.37, 0.46, 0.19, 0.31),
    @(32, 33, 294.8296768634683, 0.00277778, 0.00277778, 0.00555556, 442.2445152952025, 0.34, 0.46, 0.18, 0.38),
    @(33, 34, 349.0296502993302, 0.0138889, 0.0138889, 0.0277778, 523.5444754489954, 0.39, 0.5, 0.15, 0.4),
    @(34, 35, 364.1653337305623, 0.03055558, 0.03055558, 0.06111116, 546.2480005958435, 0.38, 0.43, 0.11, 0.37),
    @(35, 36, 381.7906987705871, 0.05833338, 0.05833338, 0.11666676, 572.6860481558806, 0.37, 0.49, 0.12, 0.3),
    @(36, 37, 410.3627591381613, 0.07777784, 0.07777784, 0.15555568, 615.5441387072419, 0.4, 0.46, 0.19, 0.33),
    @(37, 38, 408.4460994869478, 0.07222228, 0.07222228, 0.14444456, 612.6691492304217, 0.33, 0.45, 0.11, 0.32),
    @(38, 39, 377.9573794681603, 0.0555556, 0.0555556, 0.1111112, 566.9360692022404, 0.37, 0.48, 0.14, 0.33),
    @(39, 40, 354.1473291618483, 0.01111112, 0.01111112, 0.02222224, 531.2209937427724, 0.31, 0.44, 0.18, 0.32),
    @(40, 41, 356.5382138814032, 0, 0, 0, 534.8073208221048, 0.33, 0.5, 0.12, 0.36),
    @(41, 42, 444.5860015701382, 0, 0, 0, 666.8790023552073, 0.3, 0.45, 0.15, 0.35),
    @(42, 43, 515.7790399549028, 0, 0, 0, 773.6685599323542, 0.4, 0.48, 0.19, 0.4),
    @(43, 44, 517.557383961184, 0, 0, 0, 776.336075941776, 0.35, 0.43, 0.19, 0.39),
    @(44, 45, 423.5225047846375, 0, 0, 0, 635.2837571769562, 0.38, 0.46, 0.11, 0.31),
)

for ($i = 0; $i -lt $rowData.Count; $i++) {
    $r = $i + 2
    $vals = $rowData[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $vals[$j]
    }
}
